$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 363.63635
$ws.Range("I2").Value = 160
$ws.Range("K2").Value = 160
$ws.Range("M2").Value = -47
$ws.Range("H9").Value = 216.71428
$ws.Range("I9").Value = 78
$ws.Range("K9").Value = 78
$ws.Range("M9").Value = 91
$ws.Range("H11").Value = 311.16666
$ws.Range("I11").Value = 311.16666
$ws.Range("K11").Value = 311.16666
$ws.Range("M11").Value = -171.16666
$ws.Range("H17").Value = 1728.4762
$ws.Range("J17").Value = 1728.4762
$ws.Range("L17").Value = 5185.4286
$ws.Range("N17").Value = -5521.4286
$ws.Range("H88").Value = 4149.4375
$ws.Range("J88").Value = 3683.3635
$ws.Range("L88").Value = 3683.3635
$ws.Range("N88").Value = -4495.363499999999
$ws.Range("H91").Value = 4149.4375
$ws.Range("J91").Value = 3683.3635
$ws.Range("L91").Value = 3683.3635
$ws.Range("N91").Value = -6491.363499999999
$ws.Range("H96").Value = 1252.5186
$ws.Range("I96").Value = 1150.6316
$ws.Range("J96").Value = 1494.5
$ws.Range("K96").Value = 3451.8948
$ws.Range("L96").Value = 4483.5
$ws.Range("M96").Value = -2078.8948
$ws.Range("N96").Value = -7229.5
$ws.Range("H98").Value = 2979.205
$ws.Range("I98").Value = 2663.6858
$ws.Range("K98").Value = 2663.6858
$ws.Range("M98").Value = -1165.6858
$ws.Range("H122").Value = 2979.205
$ws.Range("I122").Value = 2663.6858
$ws.Range("K122").Value = 7991.057400000001
$ws.Range("M122").Value = -5541.057400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 596.64514
$ws.Range("I97").Value = 533.08
$ws.Range("K97").Value = 533.08
$ws.Range("M97").Value = -37.08000000000004
$ws.Range("H102").Value = 5076.5
$ws.Range("J102").Value = 8000
$ws.Range("L102").Value = 8000
$ws.Range("N102").Value = -11244
$ws.Range("H130").Value = 56378.2
$ws.Range("J130").Value = 66875.25
$ws.Range("L130").Value = 66875.25
$ws.Range("N130").Value = -76915.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6962423
$ws.Range("I20").Value = 13893979
$ws.Range("J20").Value = 30867.584
$ws.Range("K20").Value = 13893979
$ws.Range("L20").Value = 30867.584
$ws.Range("M20").Value = -13893732
$ws.Range("N20").Value = -31361.584
$ws.Range("H105").Value = 29412650
$ws.Range("I105").Value = 34483652
$ws.Range("K105").Value = 34483652
$ws.Range("M105").Value = -34481905
$ws.Range("H134").Value = 31035.5
$ws.Range("I134").Value = 34349.098
$ws.Range("K134").Value = 103047.294
$ws.Range("M134").Value = -100512.294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9645.6
$ws.Range("I16").Value = 6586.364
$ws.Range("K16").Value = 6586.364
$ws.Range("M16").Value = -6299.364
$ws.Range("H31").Value = 14580.2
$ws.Range("J31").Value = 21928.625
$ws.Range("L31").Value = 21928.625
$ws.Range("N31").Value = -22518.625
$ws.Range("H34").Value = 14580.2
$ws.Range("J34").Value = 21928.625
$ws.Range("L34").Value = 21928.625
$ws.Range("N34").Value = -22332.625
$ws.Range("H58").Value = 10863.979
$ws.Range("I58").Value = 5190.6895
$ws.Range("K58").Value = 5190.6895
$ws.Range("M58").Value = -4987.6895
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -18740
$ws.Range("H113").Value = 9645.6
$ws.Range("I113").Value = 6586.364
$ws.Range("K113").Value = 6586.364
$ws.Range("M113").Value = -4416.364
$ws.Range("H132").Value = 4496.6
$ws.Range("I132").Value = 4364.5454
$ws.Range("J132").Value = 4859.75
$ws.Range("K132").Value = 13093.6362
$ws.Range("L132").Value = 14579.25
$ws.Range("M132").Value = -10563.6362
$ws.Range("N132").Value = -19639.25
$ws.Range("H134").Value = 27032270
$ws.Range("I134").Value = 1297.1786
$ws.Range("K134").Value = 3891.5358
$ws.Range("M134").Value = -1356.5358
$ws.Range("H136").Value = 10863.979
$ws.Range("I136").Value = 5190.6895
$ws.Range("K136").Value = 15572.0685
$ws.Range("M136").Value = -13022.0685
$ws.Range("H141").Value = 247919.2
$ws.Range("J141").Value = 247919.2
$ws.Range("L141").Value = 247919.2
$ws.Range("N141").Value = -258279.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 36.27778
$ws.Range("I12").Value = 56.142857
$ws.Range("J12").Value = 23.636364
$ws.Range("K12").Value = 168.428571
$ws.Range("L12").Value = 70.909092
$ws.Range("M12").Value = 4.571428999999995
$ws.Range("N12").Value = -416.909092
$ws.Range("H80").Value = 22249.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 22249.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 66748.5
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -68620.5
$ws.Range("H83").Value = 22249.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 22249.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 200245.5
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -209605.5
$ws.Range("H134").Value = 6085.0835
$ws.Range("I134").Value = 1300
$ws.Range("J134").Value = 6520.091
$ws.Range("K134").Value = 3900
$ws.Range("L134").Value = 19560.273
$ws.Range("M134").Value = 1170
$ws.Range("N134").Value = -29700.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20357.416
$ws.Range("I80").Value = 17161.375
$ws.Range("J80").Value = 26749.5
$ws.Range("K80").Value = 17161.375
$ws.Range("L80").Value = 26749.5
$ws.Range("M80").Value = -16163.375
$ws.Range("N80").Value = -28745.5
$ws.Range("H83").Value = 20357.416
$ws.Range("I83").Value = 17161.375
$ws.Range("J83").Value = 26749.5
$ws.Range("K83").Value = 85806.875
$ws.Range("L83").Value = 133747.5
$ws.Range("M83").Value = -80814.875
$ws.Range("N83").Value = -143731.5
$ws.Range("H113").Value = 3378.5
$ws.Range("I113").Value = 2599.2
$ws.Range("K113").Value = 2599.2
$ws.Range("M113").Value = -429.1999999999998
$ws.Range("H132").Value = 3016.2354
$ws.Range("I132").Value = 2750.8
$ws.Range("K132").Value = 8252.400000000001
$ws.Range("M132").Value = -5722.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1639559.2
$ws.Range("I40").Value = 4491.727
$ws.Range("K40").Value = 4491.727
$ws.Range("M40").Value = -4355.727
$ws.Range("H46").Value = 2114.0833
$ws.Range("I46").Value = 725
$ws.Range("K46").Value = 725
$ws.Range("M46").Value = -537
$ws.Range("H100").Value = 3722
$ws.Range("J100").Value = 3249.75
$ws.Range("L100").Value = 3249.75
$ws.Range("N100").Value = -4331.75
$ws.Range("H122").Value = 26625784
$ws.Range("J122").Value = 3131111.2
$ws.Range("L122").Value = 9393333.600000001
$ws.Range("N122").Value = -9398233.600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 119999
$ws.Range("J56").Value = 119999
$ws.Range("L56").Value = 119999
$ws.Range("N56").Value = -121427
